{"js": "// Office.js (Word JavaScript API) edit script.\n// Applies four targeted text corrections to the report body, matching the\n// author's commit: fixes \"hypothesis\" -> \"hypothesize\" (x2), rewrites the\n// \"After analysis we found...\" paragraph to clarify the comment-count /\n// popularity finding, and adds a comma after \"this data\" in the limitations\n// paragraph.\nconst edits = [\n  { find: \"We hypothesis that posts with less words will have a greater popularity\", replace: \"We hypothesize that posts with less words will have a greater popularity\" },\n  { find: \"Time of creation, we hypothesis, would also play a large role\", replace: \"Time of creation, we hypothesize, would also play a large role\" },\n  { find: \"After analysis we found that this was not the case, we found that the number of comments a post had didn\\u2019t consistent result in popularity but in every dataset, there was a length of comment that had a higher average popularity. The more comments a post had did not result in a higher popularity, the data sets all resulted in the most popular posts having \\u201cmedium\\u201d length, by medium we mean it fell between the having the most and least comments on a post. The data sets used were \", replace: \"After analysis we found that the higher the number of comments a post had didn\\u2019t consistently result in higher average popularity. However, we did find that in every dataset, there was a \\u201cnumber of comments\\u201d that had a higher average popularity. The most popular posts had medium \\u201cnumber of comments\\u201d, by medium we mean it fell between the having the most and least comments on a post. The data sets used were \" },\n  { find: \"Another data related issue was that the posts text body was not provided for any of the datasets. If we had this data more analysis could have been done pertaining to content and entity analysis could have also been done on posts. \", replace: \"Another data related issue was that the posts text body was not provided for any of the datasets. If we had this data, more analysis could have been done pertaining to content and entity analysis could have also been done on posts. \" },\n];\n\nconst body = context.document.body;\n\nfor (const edit of edits) {\n  const results = body.search(edit.find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length !== 1) {\n    throw new Error(\n      `Expected exactly 1 match for \"${edit.find.slice(0, 40)}...\" but found ${results.items.length}`\n    );\n  }\n\n  results.items[0].insertText(edit.replace, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# Applies four targeted text corrections to the report body, matching the\n# author's commit: fixes \"hypothesis\" -> \"hypothesize\" (x2), rewrites the\n# \"After analysis we found...\" paragraph to clarify the comment-count /\n# popularity finding, and adds a comma after \"this data\" in the limitations\n# paragraph.\n\n$d = $word.ActiveDocument\n\nfunction Replace-ExactText($findText, $replaceText) {\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Text = $findText\n    $range.Find.MatchCase = $true\n    $range.Find.MatchWholeWord = $false\n    $range.Find.MatchWildcards = $false\n    $found = $range.Find.Execute()\n    if (-not $found) {\n        throw \"Replace-ExactText: text not found -> $findText\"\n    }\n    # $range now covers exactly the matched text (Find narrows it); overwrite it.\n    $range.Text = $replaceText\n}\n\n# 1) \"We hypothesis that posts...\" -> \"We hypothesize that posts...\"\nReplace-ExactText 'We hypothesis that posts with less words will have a greater popularity' 'We hypothesize that posts with less words will have a greater popularity'\n\n# 2) \"Time of creation, we hypothesis, ...\" -> \"Time of creation, we hypothesize, ...\"\nReplace-ExactText 'Time of creation, we hypothesis, would also play a large role' 'Time of creation, we hypothesize, would also play a large role'\n\n# 3) Rewrite of the \"After analysis we found...\" findings paragraph.\nReplace-ExactText 'After analysis we found that this was not the case, we found that the number of comments a post had didn\u2019t consistent result in popularity but in every dataset, there was a length of comment that had a higher average popularity. The more comments a post had did not result in a higher popularity, the data sets all resulted in the most popular posts having \u201cmedium\u201d length, by medium we mean it fell between the having the most and least comments on a post. The data sets used were ' 'After analysis we found that the higher the number of comments a post had didn\u2019t consistently result in higher average popularity. However, we did find that in every dataset, there was a \u201cnumber of comments\u201d that had a higher average popularity. The most popular posts had medium \u201cnumber of comments\u201d, by medium we mean it fell between the having the most and least comments on a post. The data sets used were '\n\n# 4) \"If we had this data more analysis...\" -> \"If we had this data, more analysis...\"\nReplace-ExactText 'Another data related issue was that the posts text body was not provided for any of the datasets. If we had this data more analysis could have been done pertaining to content and entity analysis could have also been done on posts. ' 'Another data related issue was that the posts text body was not provided for any of the datasets. If we had this data, more analysis could have been done pertaining to content and entity analysis could have also been done on posts. '\n"}
